# Fixed naive component forecaster bug - Presentation state 11.02.
# Updates the forecast-error statistics (ME, MAE, MSE, RMSE, SE, N) for
# quarters Q1..Q9 (rows 2-10) and the N count for Q9 row 11 on the
# "first_eval" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q1)
$ws.Range("B2").Value = -0.1541247462703282
$ws.Range("C2").Value = 0.852839648487878
$ws.Range("D2").Value = 1.456859763062113
$ws.Range("E2").Value = 1.207004458592475
$ws.Range("F2").Value = 1.209035744782678
$ws.Range("G2").Value = 51

# Row 3 (Q2)
$ws.Range("B3").Value = 0.1946594554508169
$ws.Range("C3").Value = 0.9013317542106349
$ws.Range("D3").Value = 1.740862592224634
$ws.Range("E3").Value = 1.319417520053692
$ws.Range("F3").Value = 1.318227895439372
$ws.Range("G3").Value = 50

# Row 4 (Q3)
$ws.Range("B4").Value = 0.09587990046027713
$ws.Range("C4").Value = 0.7813449169389153
$ws.Range("D4").Value = 1.273947737936718
$ws.Range("E4").Value = 1.128692933413122
$ws.Range("F4").Value = 1.136267503977733
$ws.Range("G4").Value = 49

# Row 5 (Q4)
$ws.Range("B5").Value = 0.2010189561291209
$ws.Range("C5").Value = 0.8357504963916088
$ws.Range("D5").Value = 1.399781063783387
$ws.Range("E5").Value = 1.183123435565109
$ws.Range("F5").Value = 1.1782594200759
$ws.Range("G5").Value = 48

# Row 6 (Q5)
$ws.Range("B6").Value = 0.1134028742995598
$ws.Range("C6").Value = 0.7037266954146911
$ws.Range("D6").Value = 0.9458163126725945
$ws.Range("E6").Value = 0.9725308800611909
$ws.Range("F6").Value = 0.976338955044726
$ws.Range("G6").Value = 47

# Row 7 (Q6)
$ws.Range("B7").Value = 0.1165933518286177
$ws.Range("C7").Value = 0.7311150576301884
$ws.Range("D7").Value = 1.063932431078233
$ws.Range("E7").Value = 1.03147100350821
$ws.Range("F7").Value = 1.038617324277796
$ws.Range("G7").Value = 38

# Row 8 (Q7)
$ws.Range("B8").Value = 0.1240158675361395
$ws.Range("C8").Value = 0.6673632257439414
$ws.Range("D8").Value = 0.9510627036727184
$ws.Range("E8").Value = 0.9752244375899931
$ws.Range("F8").Value = 0.9806497622900503
$ws.Range("G8").Value = 37

# Row 9 (Q8)
$ws.Range("B9").Value = -0.06833470819847304
$ws.Range("C9").Value = 0.3931001917406148
$ws.Range("D9").Value = 0.2556107661499554
$ws.Range("E9").Value = 0.5055796338362092
$ws.Range("F9").Value = 0.5139538519170813
$ws.Range("G9").Value = 20

# Row 10 (Q9)
$ws.Range("B10").Value = -0.001481533045873319
$ws.Range("C10").Value = 0.394696832903608
$ws.Range("D10").Value = 0.2639446156350248
$ws.Range("E10").Value = 0.5137554044825463
$ws.Range("F10").Value = 0.5347313553733598
$ws.Range("G10").Value = 13

# Row 11 (Q10) - N (G11) unchanged at 5
$ws.Range("B11").Value = 0.1038035452512703
$ws.Range("C11").Value = 0.3551067294956168
$ws.Range("D11").Value = 0.1935206269752047
$ws.Range("E11").Value = 0.4399097941342119
$ws.Range("F11").Value = 0.4779454087137884
